# Updated cryptos list on Fri Nov 17 13:28:28 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must stay plain text even if it looks like a
# pure number (e.g. "0.620", "19.80"), without leaving behind a lingering
# cell style (quote-prefix / text format) on the cell.
function Set-TextValue($range, $value) {
    $range.Style = "Normal"
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "36.539.34"
$ws.Range("E2").Value = "  -1.24%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.967.94"
$ws.Range("E3").Value = "  -3.37%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.00%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "244.94"
$ws.Range("E5").Value = "  -1.49%  "

# Row 6 - XRP
Set-TextValue $ws.Range("D6") "0.620"
$ws.Range("E6").Value = "  -2.72%  "

# Row 7 - Solana
Set-TextValue $ws.Range("D7") "59.33"
$ws.Range("E7").Value = "  -6.14%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.09%  "

# Row 9 - Cardano
Set-TextValue $ws.Range("D9") "0.377"
$ws.Range("E9").Value = "  -2.33%  "

# Row 10 - OKB
Set-TextValue $ws.Range("D10") "55.92"
$ws.Range("E10").Value = "  -4.04%  "

# Row 11 - Dogecoin
Set-TextValue $ws.Range("D11") "0.0850"
$ws.Range("E11").Value = "  +6.53%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +0.15%  "

# Row 13 - Avalanche
Set-TextValue $ws.Range("D13") "22.25"
$ws.Range("E13").Value = "  -2.24%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.256.82"
$ws.Range("E15").Value = "  -3.30%  "

# Row 16 - Chainlink
Set-TextValue $ws.Range("D16") "13.63"
$ws.Range("E16").Value = "  -5.16%  "

# Row 17 - Polkadot
Set-TextValue $ws.Range("D17") "5.37"
$ws.Range("E17").Value = "  -2.76%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "1.969.19"
$ws.Range("E18").Value = "  -3.51%  "

# Row 19 - WrappedBTC
$ws.Range("D19").Value = "36.478.53"
$ws.Range("E19").Value = "  -1.18%  "

# Row 20 - ShibaInu
$ws.Range("E20").Value = "  +0.94%  "

# Row 21 - Litecoin
Set-TextValue $ws.Range("D21") "70.65"
$ws.Range("E21").Value = "  -1.97%  "

# Row 22 - now BitcoinCash (was Uniswap)
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Range("D22") "232.05"
$ws.Range("E22").Value = "  -1.73%  "

# Row 23 - now Uniswap (was BitcoinCash)
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Range("D23") "5.12"
$ws.Range("E23").Value = "  -4.75%  "

# Row 24 - Dai
$ws.Range("E24").Value = "  +0.11%  "

# Row 25 - PancakeSwap
Set-TextValue $ws.Range("D25") "2.53"
$ws.Range("E25").Value = "  +0.48%  "

# Row 26 - Toncoin
$ws.Range("E26").Value = "  -2.79%  "

# Row 27 - Cosmos
Set-TextValue $ws.Range("D27") "9.60"
$ws.Range("E27").Value = "  -1.19%  "

# Row 28 - Monero
Set-TextValue $ws.Range("D28") "165.55"
$ws.Range("E28").Value = "  +4.05%  "

# Row 29 - EthereumClassic
Set-TextValue $ws.Range("D29") "19.80"
$ws.Range("E29").Value = "  -1.96%  "

# Row 30 - Kaspa
Set-TextValue $ws.Range("D30") "0.122"
$ws.Range("E30").Value = "  -8.83%  "

# Row 31 - Stellar
$ws.Range("E31").Value = "  -1.92%  "

# Row 32 - ImmutableX
Set-TextValue $ws.Range("D32") "1.18"
$ws.Range("E32").Value = "  +0.55%  "

# Row 33 - Filecoin
Set-TextValue $ws.Range("D33") "4.80"
$ws.Range("E33").Value = "  -5.07%  "

# Row 34 - Hedera
Set-TextValue $ws.Range("D34") "0.0642"
$ws.Range("E34").Value = "  +3.90%  "

# Row 35 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D35") "4.38"
$ws.Range("E35").Value = "  -2.53%  "

# Row 36 - THORChain
Set-TextValue $ws.Range("D36") "6.20"
$ws.Range("E36").Value = "  -1.28%  "

# Row 37 - BinanceUSD
$ws.Range("E37").Value = "  -0.01%  "

# Row 38 - WEMIXToken
$ws.Range("E38").Value = "  -1.82%  "

# Row 39 - LidoDAOToken
$ws.Range("E39").Value = "  -7.54%  "

# Row 40 - RenderToken
$ws.Range("E40").Value = "  -3.94%  "

# Row 41 - Cronos
Set-TextValue $ws.Range("D41") "0.0982"
$ws.Range("E41").Value = "  -0.78%  "

# Row 42 - TrustWalletToken
$ws.Range("E42").Value = "  -3.36%  "

# Row 43 - HuobiToken
$ws.Range("E43").Value = "  -3.44%  "

# Row 44 - VeChain
Set-TextValue $ws.Range("D44") "0.0211"

# Row 45 - InjectiveProtocol
Set-TextValue $ws.Range("D45") "15.98"
$ws.Range("E45").Value = "  -6.23%  "

# Row 46 - ARBITRUM
$ws.Range("E46").Value = "  -6.69%  "

# Row 47 - Aave
Set-TextValue $ws.Range("D47") "89.55"
$ws.Range("E47").Value = "  -4.09%  "

# Row 48 - now FraxShare (was Maker)
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D48") "7.44"
$ws.Range("E48").Value = "  -3.27%  "

# Row 49 - now Maker (was FraxShare)
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "1.353.83"
$ws.Range("E49").Value = "  -0.82%  "

# Row 50 - MXToken
$ws.Range("E50").Value = "  -3.16%  "

# Row 51 - MultiversX
Set-TextValue $ws.Range("D51") "45.32"
$ws.Range("E51").Value = "  -0.08%  "
